$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for columns C and D (they mirror the same value)
$rows = @{
    3  = 187
    5  = 949
    7  = 23
    9  = 174
    11 = 201
    13 = 4
    15 = 18
    17 = 14
    19 = 40
    21 = 5
    23 = 457
    24 = 330
    26 = 56
    28 = 333
    30 = 34
    32 = 62
    34 = 7
    36 = 69
    38 = 50
}

foreach ($r in $rows.Keys) {
    $val = $rows[$r]
    $ws.Cells.Item($r, 3).Value = $val
    $ws.Cells.Item($r, 4).Value = $val
}

# Row 39 only has column C updated (D39 does not exist)
$ws.Cells.Item(39, 3).Value = 158.578947368421
